$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.056.70"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.091.01"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.33"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.00"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.092.16"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.37"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("E12").Value = "  +3.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.625.08"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.198.93"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.087.82"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.51"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.94"
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.15"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.499"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  -6.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.80"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("E35").Value = "  +9.16%  "
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.52"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.03"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.15"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.57"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.695"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.390.36"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.62"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.132.70"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.964"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("E51").Value = "  -1.22%  "
